# Books-Stock-Data.xlsx: rework the "itemloc" sheet so that the
# "max capacity" column (originally column F) sits right after
# "qty available" (column C), pushing "expiry date" and "fifo date"
# one column to the right (D/E -> E/F). This is the classic
# "select column, Cut, select destination column, Insert Cut Cells"
# move in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("itemloc")

# Switch to the itemloc sheet (becomes the active/selected tab).
$ws.Activate()

# Cut column F ("max capacity") and insert it before column D,
# shifting the old D ("expiry date") and E ("fifo date") right by one.
$ws.Columns("F").Cut()
$ws.Columns("D").Insert()

# Leave the selection where the user ended up after the edit.
$ws.Range("J10").Select()
